# Update scripts with new TPM values.
# Re-labels the Sending/Target cluster columns (A-D) and refreshes the
# downstream expression / specificity metrics (E-T) for rows 2-7 of the
# LR-pairs (Nxph3-Nrxn1) sheet to reflect the new TPM-based computation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ordered row -> full row (A:T) values, as produced by the updated TPM run.
$rowsInOrder = @(2, 3, 4, 5, 6, 7)

$data = @{
    2 = @("FAPs", "Nxph3", "Nrxn1", "ECs",
          1, 0.3333333333333333, 0.2860933333333334, 0.85828, 0.5596515903460949, 0.5596515903460949,
          1, 0.3333333333333333, 0.002858666666666667, 0.008576, 0.002669495535069502, 0.002669495535069501,
          0.0008178454755555558, 0.00736060928, 0.001493987421623446, 0.001493987421623446)
    3 = @("FAPs", "Nxph3", "Nrxn1", "FAPs",
          1, 0.3333333333333333, 0.2860933333333334, 0.85828, 0.5596515903460949, 0.5596515903460949,
          2, 0.6666666666666666, 0.07823633333333334, 0.234709, 0.07305907503971872, 0.0730590750397187,
          0.02238289339111112, 0.20144604052, 0.04088762753519327, 0.04088762753519326)
    4 = @("FAPs", "Nxph3", "Nrxn1", "MuSCs",
          1, 0.3333333333333333, 0.2860933333333334, 0.85828, 0.5596515903460949, 0.5596515903460949,
          3, 1, 0.989769, 2.969307, 0.9242714294252118, 0.9242714294252118,
          0.28316631244, 2.54849681196, 0.5172699753892782, 0.5172699753892782)
    5 = @("MuSCs", "Nxph3", "Nrxn1", "ECs",
          1, 0.3333333333333333, 0.2251056666666666, 0.675317, 0.4403484096539051, 0.4403484096539051,
          1, 0.3333333333333333, 0.002858666666666667, 0.008576, 0.002669495535069502, 0.002669495535069501,
          0.0006435020657777778, 0.005791518592, 0.001175508113446056, 0.001175508113446056)
    6 = @("MuSCs", "Nxph3", "Nrxn1", "FAPs",
          1, 0.3333333333333333, 0.2251056666666666, 0.675317, 0.4403484096539051, 0.4403484096539051,
          2, 0.6666666666666666, 0.07823633333333334, 0.234709, 0.07305907503971872, 0.0730590750397187,
          0.01761144197255556, 0.158502977753, 0.03217144750452545, 0.03217144750452545)
    7 = @("MuSCs", "Nxph3", "Nrxn1", "MuSCs",
          1, 0.3333333333333333, 0.2251056666666666, 0.675317, 0.4403484096539051, 0.4403484096539051,
          3, 1, 0.989769, 2.969307, 0.9242714294252118, 0.9242714294252118,
          0.222802610591, 2.005223495319, 0.4070014540359336, 0.4070014540359336)
}

foreach ($r in $rowsInOrder) {
    $rowVals = $data[$r]
    for ($col = 1; $col -le $rowVals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowVals[$col - 1]
    }
}
